$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Z1 value from 2 to 3
$ws.Range("Z1").Value = 3

# Add new row 2 data (stored as text, matching the source data)
$ws.Range("A2:B2").NumberFormat = "@"
$ws.Range("A2").Value = "741710024"
$ws.Range("B2").Value = "0.5"
